# 4.1.1.1b workbook update:
#  - refresh the Kyrgyz title wording in A1
#  - add the new "2023" data column (K), copying the formatting of column J
#  - apply vertical-center alignment across the used range
#  - nudge a few row heights to match the refreshed layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Kyrgyz title (A1) to the new translation ---
$ws.Range("A1").Value = "4.1.1.1b Башталгыч билим берүүнү аяктоо деңгээли"

# --- Add the new "2023" column (K), matching the formatting of the preceding column (J) ---
$ws.Range("J3:J14").Copy()
$ws.Range("K3:K14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("K4").Value = 2023

$ws.Range("K5").Value = 100.4062134821358
$ws.Range("K6").Value = 104.06287706751064
$ws.Range("K7").Value = 102.93580727245744
$ws.Range("K8").Value = 99.942837442404183
$ws.Range("K9").Value = 106.42286904738506
$ws.Range("K10").Value = 103.81318629146574
$ws.Range("K11").Value = 104.42396096858288
$ws.Range("K12").Value = 101.98873952645914
$ws.Range("K13").Value = 90.624818310428424
$ws.Range("K14").Value = 93.159844818577312

# --- Apply vertical-center alignment across the whole used range (A1:K14) ---
$ws.Range("A1:K14").VerticalAlignment = -4108

# --- Minor row-height adjustments to match the refreshed layout ---
$ws.Rows.Item(2).RowHeight = 14.25
foreach ($r in 4..14) {
    $ws.Rows.Item($r).RowHeight = 13.5
}

# --- Reset the active selection back to the top-left cell ---
$null = $ws.Range("A1").Select()
